# Fix "Excel file total marks error"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right count 5 -> 4, Wrong marking -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Total right marks 95 -> 76, wrong marks -1 -> -2
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -2

# Update the displayed total string to match the corrected score
$ws.Range("E12").Value = "74 / 112"
